$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-271) stores a date serial number that needs to move
# from 45202 (2023-10-03) to 45203 (2023-10-04).
$ws.Range("C2:C271").Value = 45203
